$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "353.70", "1.00")
# that Excel would otherwise auto-convert to a Number and mangle (dropping
# trailing zeros). Force text via NumberFormat "@" while assigning, then
# ClearFormats() so the cell keeps the General style of the original file
# (only the stored value/type changes, not the cell style).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.815.38"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.805.78"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.30"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.94"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.248.99"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.809.74"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.821.45"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +8.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.18"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.06"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +14.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.39"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0455"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0894"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.06%  "
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.72%  "
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +10.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.106.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.964"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("E51").Value = "  +8.60%  "
